$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3460
$ws.Range("I18").Value = 3460
$ws.Range("K18").Value = 3460
$ws.Range("M18").Value = -3176
$ws.Range("H55").Value = 163.84616
$ws.Range("I55").Value = 170
$ws.Range("J55").Value = 162
$ws.Range("K55").Value = 170
$ws.Range("L55").Value = 162
$ws.Range("M55").Value = 44
$ws.Range("N55").Value = -590
$ws.Range("H86").Value = 2755.3333
$ws.Range("I86").Value = 2685.5715
$ws.Range("K86").Value = 2685.5715
$ws.Range("M86").Value = -1562.5715
$ws.Range("H89").Value = 2755.3333
$ws.Range("I89").Value = 2685.5715
$ws.Range("K89").Value = 13427.8575
$ws.Range("M89").Value = -7811.8575
$ws.Range("H137").Value = 9404.474
$ws.Range("I137").Value = 1565.9166
$ws.Range("J137").Value = 22842
$ws.Range("K137").Value = 4697.7498
$ws.Range("L137").Value = 68526
$ws.Range("M137").Value = -2147.7498
$ws.Range("N137").Value = -73626
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2657
$ws.Range("I5").Value = 218.27272
$ws.Range("J5").Value = 5637.6665
$ws.Range("K5").Value = 218.27272
$ws.Range("L5").Value = 5637.6665
$ws.Range("M5").Value = -106.27272
$ws.Range("N5").Value = -5861.6665
$ws.Range("H32").Value = 302130.5
$ws.Range("I32").Value = 340747.53
$ws.Range("K32").Value = 340747.53
$ws.Range("M32").Value = -340460.53
$ws.Range("H45").Value = 2950.7144
$ws.Range("I45").Value = 2231
$ws.Range("K45").Value = 2231
$ws.Range("M45").Value = -1854
$ws.Range("H61").Value = 10189.556
$ws.Range("I61").Value = 2266.3333
$ws.Range("J61").Value = 14151.167
$ws.Range("K61").Value = 2266.3333
$ws.Range("L61").Value = 14151.167
$ws.Range("M61").Value = -2054.3333
$ws.Range("N61").Value = -14575.167
$ws.Range("H74").Value = 5453.396
$ws.Range("I74").Value = 945.4211
$ws.Range("K74").Value = 945.4211
$ws.Range("M74").Value = -71.42110000000002
$ws.Range("H77").Value = 5453.396
$ws.Range("I77").Value = 945.4211
$ws.Range("K77").Value = 4727.1055
$ws.Range("M77").Value = -359.1054999999997
$ws.Range("H122").Value = 3122.4285
$ws.Range("I122").Value = 2971.6
$ws.Range("K122").Value = 8914.799999999999
$ws.Range("M122").Value = -6464.799999999999
$ws.Range("H132").Value = 1045438.75
$ws.Range("I132").Value = 1140005.9
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 3420017.7
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -3417487.7
$ws.Range("N132").Value = -20660
$ws.Range("H134").Value = 64992
$ws.Range("J134").Value = 64992
$ws.Range("L134").Value = 64992
$ws.Range("N134").Value = -75132
$ws.Range("H136").Value = 10189.556
$ws.Range("I136").Value = 2266.3333
$ws.Range("J136").Value = 14151.167
$ws.Range("K136").Value = 6798.999899999999
$ws.Range("L136").Value = 42453.501
$ws.Range("M136").Value = -4248.999899999999
$ws.Range("N136").Value = -47553.501
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2657
$ws.Range("I4").Value = 218.27272
$ws.Range("J4").Value = 5637.6665
$ws.Range("K4").Value = 218.27272
$ws.Range("L4").Value = 5637.6665
$ws.Range("M4").Value = -103.27272
$ws.Range("N4").Value = -5867.6665
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").Value = $null
$ws.Range("H134").Value = 6321.6206
$ws.Range("I134").Value = 2897.1924
$ws.Range("J134").Value = 36000
$ws.Range("K134").Value = 8691.5772
$ws.Range("L134").Value = 108000
$ws.Range("M134").Value = -6156.5772
$ws.Range("N134").Value = -113070
$ws.Range("H135").Value = 79999
$ws.Range("J135").Value = 79999
$ws.Range("L135").Value = 79999
$ws.Range("N135").Value = -90139
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("N137").Value = $null
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").Value = $null
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5899.273
$ws.Range("J94").Value = 1377.5
$ws.Range("L94").Value = 1377.5
$ws.Range("N94").Value = -2279.5
$ws.Range("H116").Value = 74000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = $null
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("N117").Value = $null
$ws.Range("H132").Value = 2333
$ws.Range("I132").Value = 1912
$ws.Range("K132").Value = 5736
$ws.Range("M132").Value = -3206
$ws.Range("H134").Value = 4331
$ws.Range("I134").Value = 4331
$ws.Range("K134").Value = 12993
$ws.Range("M134").Value = -10458
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2871.647
$ws.Range("J131").Value = 2871.647
$ws.Range("L131").Value = 8614.940999999999
$ws.Range("N131").Value = -18694.941
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 51667.36
$ws.Range("I70").Value = 42247.633
$ws.Range("K70").Value = 42247.633
$ws.Range("M70").Value = -41977.633
$ws.Range("H73").Value = 51667.36
$ws.Range("I73").Value = 42247.633
$ws.Range("K73").Value = 42247.633
$ws.Range("M73").Value = -41311.633
$ws.Range("H97").Value = 784.06976
$ws.Range("I97").Value = 796.1622
$ws.Range("J97").Value = 709.5
$ws.Range("K97").Value = 796.1622
$ws.Range("L97").Value = 709.5
$ws.Range("M97").Value = -300.1622
$ws.Range("N97").Value = -1701.5
$ws.Range("H132").Value = 10477.403
$ws.Range("I132").Value = 12219.728
$ws.Range("J132").Value = 4580.3076
$ws.Range("K132").Value = 36659.18399999999
$ws.Range("L132").Value = 13740.9228
$ws.Range("M132").Value = -34129.18399999999
$ws.Range("N132").Value = -18800.9228
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2153.577
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H27").Value = 2153.577
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = $null
$ws.Range("H55").Value = 1756.909
$ws.Range("J55").Value = 1766.9048
$ws.Range("L55").Value = 1766.9048
$ws.Range("N55").Value = -2112.9048
$ws.Range("H122").Value = 7000.4
$ws.Range("I122").Value = 3668
$ws.Range("J122").Value = 8428.571
$ws.Range("K122").Value = 11004
$ws.Range("L122").Value = 25285.713
$ws.Range("M122").Value = -8554
$ws.Range("N122").Value = -30185.713
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2350
$ws.Range("J81").Value = 3670
$ws.Range("L81").Value = 7340
$ws.Range("N81").Value = -9462
$ws.Range("H84").Value = 2350
$ws.Range("J84").Value = 3670
$ws.Range("L84").Value = 36700
$ws.Range("N84").Value = -47308
$ws.Range("H132").Value = 2776.524
$ws.Range("I132").Value = 2243.875
$ws.Range("K132").Value = 6731.625
$ws.Range("M132").Value = -4201.625
$ws.Range("H136").Value = 2178
$ws.Range("I136").Value = 2178
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6534
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null
